# "Videogame to game change"
#
# 1) The cached text of every `datetimeFigureOut` field ("today" date
#    placeholder) in the slide master and the 11 slide layouts is bumped
#    from 12/31/2020 -> 1/12/2021 (this is the auto-updating date field
#    PowerPoint re-stamps whenever the deck is saved).
# 2) Slide 1's title "2D VideoGame Development" is corrected to
#    "2D Game Development" (merging the three runs that spelled out
#    "2D " / "VideoGame" / " Development" into a single run).
#
# NOTE: this runtime's PowerShell variables are not function-scoped, so
# every loop below uses a uniquely named counter to avoid cross-loop
# interference.

$p = $ppt.ActivePresentation
$newDate = "1/12/2021"

# Slide master's own Date Placeholder.
$masterShapes = $p.SlideMaster.Shapes
for ($mi = 1; $mi -le $masterShapes.Count; $mi++) {
    $masterShape = $masterShapes.Item($mi)
    if ($masterShape.Name -like "Date Placeholder*") {
        $masterShape.TextFrame.TextRange.Text = $newDate
    }
}

# Every slide layout's Date Placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    $layoutShapes = $layout.Shapes
    for ($lj = 1; $lj -le $layoutShapes.Count; $lj++) {
        $layoutShape = $layoutShapes.Item($lj)
        if ($layoutShape.Name -like "Date Placeholder*") {
            $layoutShape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Notes master's Date Placeholder (some hosts expose this object
# read-only; attempt it but do not fail the whole script if so).
try {
    $notesShapes = $p.NotesMaster.Shapes
    for ($ni = 1; $ni -le $notesShapes.Count; $ni++) {
        $notesShape = $notesShapes.Item($ni)
        if ($notesShape.Name -like "Date Placeholder*") {
            $notesShape.TextFrame.TextRange.Text = $newDate
        }
    }
} catch {
}

# Fix the title text on slide 1 ("Title 1" placeholder).
$slide1 = $p.Slides.Item(1)
$slide1Shapes = $slide1.Shapes
for ($si = 1; $si -le $slide1Shapes.Count; $si++) {
    $titleShape = $slide1Shapes.Item($si)
    if ($titleShape.Name -eq "Title 1") {
        $titleRange = $titleShape.TextFrame.TextRange
        $titleText = $titleRange.Text
        $oldTitle = "2D VideoGame Development"
        $titleIdx = $titleText.IndexOf($oldTitle)
        if ($titleIdx -ge 0) {
            $titleSub = $titleRange.Characters($titleIdx + 1, $oldTitle.Length)
            $titleSub.Text = "2D Game Development"
        }
    }
}
